$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: header row (bold "title" style) ---
$ws.Range("B9").Value = "Number of employees"
$ws.Range("C9").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D9").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B9:D9").Font.Bold = $true

# --- Row 10: Micro ---
$ws.Range("A10").Value = "Micro"
$ws.Range("B10").Value = "1-10"
$ws.Range("C10").Value = "'"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'"
$ws.Range("D10").Style = "Normal"

# --- Row 11: Small ---
$ws.Range("A11").Value = "Small"
$ws.Range("B11").Value = "11-25"
$ws.Range("C11").Value = "'"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'"
$ws.Range("D11").Style = "Normal"

# --- Row 12: Medium ---
$ws.Range("A12").Value = "Medium"
$ws.Range("B12").Value = "26-100"
$ws.Range("C12").Value = "'"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'"
$ws.Range("D12").Style = "Normal"

# --- Row 13: Large ---
$ws.Range("A13").Value = "Large"
$ws.Range("B13").Value = ">100"
$ws.Range("C13").Value = "'"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'"
$ws.Range("D13").Style = "Normal"
